# Applies the scheduled-runner price/profit refresh to the Excalibur_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 2750
$ws.Range("J21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("N21").Value = -5936
$ws.Range("H23").Value = 2750
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5468
$ws.Range("H86").Value = 4800.0435
$ws.Range("I86").Value = 3920.111
$ws.Range("K86").Value = 3920.111
$ws.Range("M86").Value = -2797.111
$ws.Range("H89").Value = 4800.0435
$ws.Range("I89").Value = 3920.111
$ws.Range("K89").Value = 19600.555
$ws.Range("M89").Value = -13984.555
$ws.Range("H92").Value = 972.1905
$ws.Range("I92").Value = 542.7059
$ws.Range("J92").Value = 2797.5
$ws.Range("K92").Value = 542.7059
$ws.Range("L92").Value = 2797.5
$ws.Range("M92").Value = 705.2941
$ws.Range("N92").Value = -5293.5
$ws.Range("H112").Value = 5067.0464
$ws.Range("J112").Value = 5320.623
$ws.Range("L112").Value = 15961.869
$ws.Range("N112").Value = -18177.869
$ws.Range("H129").Value = 1661.1333
$ws.Range("I129").Value = 1405.5834
$ws.Range("K129").Value = 4216.7502
$ws.Range("M129").Value = 783.2497999999996
$ws.Range("H132").Value = 66773.61
$ws.Range("I132").Value = 76746.94
$ws.Range("K132").Value = 230240.82
$ws.Range("M132").Value = -227710.82
$ws.Range("H137").Value = 1207379.6
$ws.Range("I137").Value = 2901.4167
$ws.Range("K137").Value = 8704.250100000001
$ws.Range("M137").Value = -6154.250100000001
$ws.Range("H138").Value = 3343.9456
$ws.Range("I138").Value = 1502.7391
$ws.Range("J138").Value = 4667.3125
$ws.Range("K138").Value = 4508.2173
$ws.Range("L138").Value = 14001.9375
$ws.Range("M138").Value = 631.7826999999997
$ws.Range("N138").Value = -24281.9375
$ws.Range("H141").Value = 2620.6667
$ws.Range("I141").Value = 2620.6667
$ws.Range("K141").Value = 7862.000100000001
$ws.Range("M141").Value = -2682.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23811.266
$ws.Range("I32").Value = 24203.188
$ws.Range("K32").Value = 24203.188
$ws.Range("M32").Value = -23916.188
$ws.Range("H75").Value = 74173
$ws.Range("J75").Value = 74173
$ws.Range("L75").Value = 74173
$ws.Range("N75").Value = -75921
$ws.Range("H78").Value = 74173
$ws.Range("J78").Value = 74173
$ws.Range("L78").Value = 222519
$ws.Range("N78").Value = -231255
$ws.Range("H97").Value = 1268.2273
$ws.Range("I97").Value = 842.6429000000001
$ws.Range("J97").Value = 2013
$ws.Range("K97").Value = 842.6429000000001
$ws.Range("L97").Value = 2013
$ws.Range("M97").Value = -346.6429000000001
$ws.Range("N97").Value = -3005
$ws.Range("H132").Value = 627097.4399999999
$ws.Range("I132").Value = 770793
$ws.Range("K132").Value = 2312379
$ws.Range("M132").Value = -2309849

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1546.5
$ws.Range("I94").Value = 1544.8182
$ws.Range("K94").Value = 1544.8182
$ws.Range("M94").Value = -1093.8182
$ws.Range("H105").Value = 2377.2222
$ws.Range("I105").Value = 2377.2222
$ws.Range("K105").Value = 2377.2222
$ws.Range("M105").Value = -630.2222000000002
$ws.Range("H134").Value = 3155604.5
$ws.Range("I134").Value = 3404512
$ws.Range("J134").Value = 2284428
$ws.Range("K134").Value = 10213536
$ws.Range("L134").Value = 6853284
$ws.Range("M134").Value = -10211001
$ws.Range("N134").Value = -6858354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28618.117
$ws.Range("I31").Value = 23402.75
$ws.Range("K31").Value = 23402.75
$ws.Range("M31").Value = -23107.75
$ws.Range("H34").Value = 28618.117
$ws.Range("I34").Value = 23402.75
$ws.Range("K34").Value = 23402.75
$ws.Range("M34").Value = -23200.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9648.777
$ws.Range("I3").Value = 9548.429
$ws.Range("K3").Value = 28645.287
$ws.Range("M3").Value = -28533.287
$ws.Range("H4").Value = 24442790
$ws.Range("I4").Value = 37000544
$ws.Range("K4").Value = 111001632
$ws.Range("M4").Value = -111001520
$ws.Range("H11").Value = 28.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 28.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 85.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -365.5
$ws.Range("H92").Value = 668
$ws.Range("J92").Value = 800
$ws.Range("L92").Value = 2400
$ws.Range("N92").Value = -4896
$ws.Range("H139").Value = 4571.8335
$ws.Range("I139").Value = 3143.6667
$ws.Range("K139").Value = 9431.000100000001
$ws.Range("M139").Value = -4291.000100000001
$ws.Range("H140").Value = 2923.4482
$ws.Range("I140").Value = 1759
$ws.Range("J140").Value = 4171.0713
$ws.Range("K140").Value = 5277
$ws.Range("L140").Value = 12513.2139
$ws.Range("M140").Value = -97
$ws.Range("N140").Value = -22873.2139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1344021.9
$ws.Range("I132").Value = 1510274.6
$ws.Range("K132").Value = 4530823.800000001
$ws.Range("M132").Value = -4528293.800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1366
$ws.Range("J22").Value = 2375
$ws.Range("L22").Value = 2375
$ws.Range("N22").Value = -2965
$ws.Range("H27").Value = 1366
$ws.Range("J27").Value = 2375
$ws.Range("L27").Value = 2375
$ws.Range("N27").Value = -2589
$ws.Range("H40").Value = 3005.2632
$ws.Range("I40").Value = 3136.9333
$ws.Range("J40").Value = 2511.5
$ws.Range("K40").Value = 3136.9333
$ws.Range("L40").Value = 2511.5
$ws.Range("M40").Value = -3000.9333
$ws.Range("N40").Value = -2783.5
$ws.Range("H76").Value = 11858.667
$ws.Range("J76").Value = 11858.667
$ws.Range("L76").Value = 11858.667
$ws.Range("N76").Value = -12534.667
$ws.Range("H79").Value = 11858.667
$ws.Range("J79").Value = 11858.667
$ws.Range("L79").Value = 11858.667
$ws.Range("N79").Value = -14198.667
$ws.Range("H130").Value = 90357.336
$ws.Range("J130").Value = 90357.336
$ws.Range("L130").Value = 90357.336
$ws.Range("N130").Value = -100397.336
$ws.Range("H132").Value = 943685.9399999999
$ws.Range("I132").Value = 1124979.9
$ws.Range("J132").Value = 7000.3335
$ws.Range("K132").Value = 3374939.7
$ws.Range("L132").Value = 21001.0005
$ws.Range("M132").Value = -3372409.7
$ws.Range("N132").Value = -26061.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3006022
$ws.Range("I132").Value = 3196573.8
$ws.Range("J132").Value = 4830
$ws.Range("K132").Value = 9589721.399999999
$ws.Range("L132").Value = 14490
$ws.Range("M132").Value = -9587191.399999999
$ws.Range("N132").Value = -19550
